# Generate Report for handoff
#
# The source file "9885c67c-...md" finished its handoff cycle and was
# replaced by a new source file "c8ad2577-...md" (new handoff hashes /
# timestamps), and a second source file "b89b8d5d-...md" was picked up
# whose handoff transform failed. The old "row 3" (.localization-config)
# slides down to row 4 on every sheet to make room for the new
# "Handoff transform failed" row.

$wb = $excel.ActiveWorkbook

$repoBase   = "https://github.com/OpenLocalizationTest/oltest/blob/df377a7694531c355d21530b805286fefdc45381"
$zhHoBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a10b74c2bdd01352fb4d2b55dc7aad18eacd37fd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht"
$deHoBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7d33e1ec63ca827ded427230f55aea192bccf096/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht"

$mdFile1    = "c8ad2577-c386-47d7-91e1-91f504cd81b0.md"
$mdFile2    = "b89b8d5d-93b9-40fb-b029-edd733fe590c.md"
$cfgFile    = ".localization-config"

$zhXlf      = "c8ad2577-c386-47d7-91e1-91f504cd81b0.dbe3e44bfb41595194553b8871d6be6f580bdda9.zh-cn.xlf"
$deXlf      = "c8ad2577-c386-47d7-91e1-91f504cd81b0.dbe3e44bfb41595194553b8871d6be6f580bdda9.de-de.xlf"

$zhHandoffDatetime = "2016-02-17 03:15:48"
$deHandoffDatetime = "2016-02-17 03:15:58"
$epoch             = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

$ws1.Range("B3").Value = "Handoff transform failed"
$ws1.Range("C3").Value = "Handoff transform failed"

$ws1.Range("B4").Value = "Not to be localized"
$ws1.Range("C4").Value = "Not to be localized"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "$repoBase/e2e/$mdFile1", $null, $null, $mdFile1)
$ws1.Hyperlinks.Add($ws1.Range("A3"), "$repoBase/e2e/$mdFile2", $null, $null, $mdFile2)
$ws1.Hyperlinks.Add($ws1.Range("A4"), "$repoBase/$cfgFile", $null, $null, $cfgFile)

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("D2").Value = $zhHandoffDatetime
$ws2.Range("G2").Value = $epoch
$ws2.Range("H2").Value = "Include"

$ws2.Range("B3").Value = "Handoff transform failed"
$ws2.Range("D3").Value = $epoch
$ws2.Range("G3").Value = $epoch
$ws2.Range("H3").Value = "Ignored"

$ws2.Range("B4").Value = "Not to be localized"
$ws2.Range("D4").Value = $epoch
$ws2.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("G4").Value = $epoch
$ws2.Range("H4").Value = "Ignored"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "$repoBase/e2e/$mdFile1", $null, $null, $mdFile1)
$ws2.Hyperlinks.Add($ws2.Range("C2"), "$zhHoBase/$zhXlf", $null, $null, $zhXlf)
$ws2.Hyperlinks.Add($ws2.Range("A3"), "$repoBase/e2e/$mdFile2", $null, $null, $mdFile2)
$ws2.Hyperlinks.Add($ws2.Range("A4"), "$repoBase/$cfgFile", $null, $null, $cfgFile)

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("D2").Value = $deHandoffDatetime
$ws3.Range("G2").Value = $epoch
$ws3.Range("H2").Value = "Include"

$ws3.Range("B3").Value = "Handoff transform failed"
$ws3.Range("D3").Value = $epoch
$ws3.Range("G3").Value = $epoch
$ws3.Range("H3").Value = "Ignored"

$ws3.Range("B4").Value = "Not to be localized"
$ws3.Range("D4").Value = $epoch
$ws3.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("G4").Value = $epoch
$ws3.Range("H4").Value = "Ignored"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "$repoBase/e2e/$mdFile1", $null, $null, $mdFile1)
$ws3.Hyperlinks.Add($ws3.Range("C2"), "$deHoBase/$deXlf", $null, $null, $deXlf)
$ws3.Hyperlinks.Add($ws3.Range("A3"), "$repoBase/e2e/$mdFile2", $null, $null, $mdFile2)
$ws3.Hyperlinks.Add($ws3.Range("A4"), "$repoBase/$cfgFile", $null, $null, $cfgFile)
